$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Bus sheet: widen reactive power limits and add active/reactive power
# to bus 2 so GFM inverter has an infinite admittance reference to
# compare against.
# ---------------------------------------------------------------------
$wsBus = $wb.Worksheets.Item("Bus")
$wsBus.Range("I5").Value = -999
$wsBus.Range("J5").Value = 999

$wsBus.Range("E6").Value = 0.5
$wsBus.Range("F6").Value = 0
$wsBus.Range("I6").Value = -999
$wsBus.Range("J6").Value = 999

$wsBus.Activate()
$wsBus.Range("F7").Select()

# ---------------------------------------------------------------------
# Apparatus sheet: insert a new sub-header row describing the GFM
# inverter parameters, then fill in the actual values/formulas for the
# GFM apparatus (row 6 after the insert).
# ---------------------------------------------------------------------
$wsApp = $wb.Worksheets.Item("Apparatus")

$wsApp.Rows("3:3").Insert()

$wsApp.Range("C3").Value = "wLf (pu)"
$wsApp.Range("D3").Value = "Rf (pu)"
$wsApp.Range("E3").Value = "wCf (pu)"
$wsApp.Range("F3").Value = "wLc (pu)"
$wsApp.Range("G3").Value = "Rc (pu)"
$wsApp.Range("H3").Value = "Xov (pu)"
$wsApp.Range("I3").Value = "Droop Dw"
$wsApp.Range("J3").Value = "BW droop (Hz)"
$wsApp.Range("K3").Value = "BW vdq (Hz)"
$wsApp.Range("L3").Value = "BW idq (Hz)"

$wsApp.Range("D6").Formula = "=C6/10"
$wsApp.Range("E6").Value = 0.02
$wsApp.Range("F6").Value = 0.01
$wsApp.Range("G6").Formula = "=F6/10"
$wsApp.Range("H6").Value = 0
$wsApp.Range("I6").Value = 0.1
$wsApp.Range("J6").Value = 0.5
$wsApp.Range("K6").Value = 300
$wsApp.Range("L6").Value = 600

# Narrow the newly-repurposed columns (they used to hold long labels
# inline, now the labels live in row 3 so the columns can be narrower).
$wsApp.Columns.Item(6).ColumnWidth = 66/7
$wsApp.Columns.Item(7).ColumnWidth = 87/7
$wsApp.Columns.Item(8).ColumnWidth = 81/7
$wsApp.Columns.Item(9).ColumnWidth = 82/7
$wsApp.Columns.Item(10).ColumnWidth = 100/7

$wsApp.Activate()
$wsApp.Range("E6").Select()

# ---------------------------------------------------------------------
# NetworkLine_IEEE sheet: cosmetic selection update only.
# ---------------------------------------------------------------------
$wsIeee = $wb.Worksheets.Item("NetworkLine_IEEE")
$wsIeee.Activate()
$wsIeee.Range("E9").Select()

# ---------------------------------------------------------------------
# NetworkLine sheet: change the R (pu) formula denominator from 5 to 10
# to stay consistent with the Apparatus sheet's Rf/Rc formulas.
# ---------------------------------------------------------------------
$wsLine = $wb.Worksheets.Item("NetworkLine")
$wsLine.Range("C4").Formula = "=D4/10"

$wsLine.Activate()
$wsLine.Range("C5").Select()

# ---------------------------------------------------------------------
# Advance sheet: disable the Simulink model auto-open, and update the
# selection.
# ---------------------------------------------------------------------
$wsAdv = $wb.Worksheets.Item("Advance")
$wsAdv.Range("B8").Value = 0

$wsAdv.Activate()
$wsAdv.Range("B9").Select()

# ---------------------------------------------------------------------
# Restore Apparatus as the active tab (matches workbookView activeTab).
# ---------------------------------------------------------------------
$wsApp.Activate()
